# Add a "Price" column (D) to Table2 and populate it, per the commit:
# "Added new DecimalPlaces property for ClassToExcelRowAttribute -
#  ClassToExcelRowConverter will now write to private properties."
# This re-generated the sample workbook with an extra numeric column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow Table2 from A1:C7 to A1:D7 - this adds a "Column4" header in D1
# and keeps the rest of the table structure (style, autofilter, etc.) intact.
$lo2 = $ws.ListObjects.Item("Table2")
$lo2.Resize($ws.Range("A1:D7"))

# Header
$ws.Range("D1").Value = "Price"

# Row 2: Beer
$ws.Range("D2").Value = 2.154

# Row 3: Wine (previously missing - blank row 3 in the original sheet)
$ws.Range("A3").Value = "Wine"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 10.257

# Row 4: Pepsi
$ws.Range("A4").Value = "Pepsi"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1.25

# Row 5: Coke
$ws.Range("A5").Value = "Coke"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 1.26

# Row 6: Dr. Pepper
$ws.Range("A6").Value = "Dr. Pepper"
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 1.27

# Row 7: now holds the "Avg. Number of Liters" figure instead of Dr. Pepper's
# old quantity row.
$ws.Range("A7").Value = "Avg. Number of Liters"
$ws.Range("C7").Value = 0.789
$ws.Range("D7").ClearContents()

# Reflect the new selection left behind by the edit.
$ws.Range("A8").Select()
